# Update latest output (run 168)
# - Schedule sheet: row 3 gets new aggregated values, row 4 (now redundant) is removed.
# - Detailed sheet: refreshed price/type/pump-status values for the latest run.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Schedule": update row 3, then delete row 4 (merged into row 3's window)
# ---------------------------------------------------------------------------
$wsSchedule = $wb.Worksheets.Item("Schedule")

$wsSchedule.Range("A3").Value2 = 46064.16666666666
$wsSchedule.Range("B3").Value2 = 46064.66666666666
$wsSchedule.Range("C3").Value2 = 12
$wsSchedule.Range("D3").Value2 = 45.36
$wsSchedule.Range("E3").Value2 = 1536.192411
$wsSchedule.Range("F3").Value2 = 33.86667572751323

$wsSchedule.Rows.Item(4).Delete()

# ---------------------------------------------------------------------------
# Sheet "Detailed": refresh Price / Type / Pump_Status for rows 38-97
# ---------------------------------------------------------------------------
$wsDetailed = $wb.Worksheets.Item("Detailed")

$wsDetailed.Range("B38").Value2 = 110.17505

$wsDetailed.Range("B39").Value2 = 161.25495

$wsDetailed.Range("B40").Value2 = 206.75779
$wsDetailed.Range("C40").Value2 = "historical"

$wsDetailed.Range("B41").Value2 = 299.99
$wsDetailed.Range("C41").Value2 = "historical"

$wsDetailed.Range("B42").Value2 = 299.99
$wsDetailed.Range("C42").Value2 = "historical"

$wsDetailed.Range("B43").Value2 = 222.23033
$wsDetailed.Range("C43").Value2 = "historical"

$wsDetailed.Range("B44").Value2 = 178.67423
$wsDetailed.Range("C44").Value2 = "historical"

$wsDetailed.Range("C45").Value2 = "historical"

$wsDetailed.Range("B46").Value2 = 144.62829
$wsDetailed.Range("C46").Value2 = "historical"

$wsDetailed.Range("B47").Value2 = 154.2
$wsDetailed.Range("C47").Value2 = "historical"

$wsDetailed.Range("B48").Value2 = 166.99
$wsDetailed.Range("C48").Value2 = "historical"

$wsDetailed.Range("B49").Value2 = 166.99
$wsDetailed.Range("C49").Value2 = "historical"

$wsDetailed.Range("B50").Value2 = 150.37229

$wsDetailed.Range("B51").Value2 = 138.42

$wsDetailed.Range("B52").Value2 = 134.72311

$wsDetailed.Range("B53").Value2 = 113.50454

$wsDetailed.Range("B54").Value2 = 105.79
$wsDetailed.Range("E54").Value2 = "OFF"

$wsDetailed.Range("B55").Value2 = 105.79
$wsDetailed.Range("E55").Value2 = "OFF"

$wsDetailed.Range("B56").Value2 = 105.79
$wsDetailed.Range("E56").Value2 = "OFF"

$wsDetailed.Range("B57").Value2 = 105.79
$wsDetailed.Range("E57").Value2 = "OFF"

$wsDetailed.Range("B58").Value2 = 108.89

$wsDetailed.Range("B59").Value2 = 84.79000000000001

$wsDetailed.Range("B60").Value2 = 93.77001

$wsDetailed.Range("B61").Value2 = 118.44767

$wsDetailed.Range("B62").Value2 = 122.19086
$wsDetailed.Range("E62").Value2 = "ON"

$wsDetailed.Range("B63").Value2 = 119.50455
$wsDetailed.Range("E63").Value2 = "ON"

$wsDetailed.Range("B64").Value2 = 59.79985
$wsDetailed.Range("E64").Value2 = "ON"

$wsDetailed.Range("B65").Value2 = 48.3489

$wsDetailed.Range("B66").Value2 = 36.07

$wsDetailed.Range("B67").Value2 = 48.30146

$wsDetailed.Range("B68").Value2 = 36.07

$wsDetailed.Range("B69").Value2 = 50.62162

$wsDetailed.Range("B70").Value2 = 36.07

$wsDetailed.Range("B71").Value2 = 36.07

$wsDetailed.Range("B72").Value2 = 36.07

$wsDetailed.Range("B73").Value2 = 36.07

$wsDetailed.Range("B74").Value2 = 36.07

$wsDetailed.Range("B75").Value2 = 52.88049

$wsDetailed.Range("B76").Value2 = 57.08

$wsDetailed.Range("B77").Value2 = 62.1786

$wsDetailed.Range("B78").Value2 = 65.34656

$wsDetailed.Range("B79").Value2 = 64.89

$wsDetailed.Range("B80").Value2 = 68.08553000000001

$wsDetailed.Range("B81").Value2 = 97.96586000000001
$wsDetailed.Range("E81").Value2 = "ON"

$wsDetailed.Range("B82").Value2 = 134.56236

$wsDetailed.Range("B83").Value2 = 248.88

$wsDetailed.Range("B84").Value2 = 12131.28128

$wsDetailed.Range("B85").Value2 = 12182.58465

$wsDetailed.Range("B86").Value2 = 12291.53742

$wsDetailed.Range("B87").Value2 = 13979.99364

$wsDetailed.Range("B91").Value2 = 248.88

$wsDetailed.Range("B93").Value2 = 133.03

$wsDetailed.Range("B94").Value2 = 126.46787

$wsDetailed.Range("B95").Value2 = 131.40716

$wsDetailed.Range("B96").Value2 = 90.43129999999999

$wsDetailed.Range("B97").Value2 = 81.16719999999999
